# Auto-generated edit script applying the cryptos.xlsx data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.345.85'
$ws.Range('E2').Value = '  +3.93%  '
$ws.Range('D3').Value = '2.600.55'
$ws.Range('E3').Value = '  +2.16%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '522.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.566'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.56%  '
$ws.Range('D9').Value = '2.624.89'
$ws.Range('E9').Value = '  +2.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.57'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('E11').Value = '  +1.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.332'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').Value = '3.067.56'
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').Value = '59.274.04'
$ws.Range('E15').Value = '  +3.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.45'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('D17').Value = '2.612.66'
$ws.Range('E17').Value = '  +2.94%  '
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '339.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.02%  '
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('E22').Value = '  +7.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.33%  '
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = '0.0₃0727'
$ws.Range('E30').Value = '  -3.33%  '
$ws.Range('E31').Value = '  -4.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.02'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('E35').Value = '  +1.25%  '
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.64%  '
$ws.Range('E38').Value = '  +3.76%  '
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.821'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.55'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '278.27'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.994'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.72'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('E45').Value = '  +2.43%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0522'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').Value = '1.989.10'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.62'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0221'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.23%  '

Write-Output "Applied 83 cell updates"
